$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: fill in the newly-added sub_data_type column and fix the
#     file_ending value (patient obj files are now "cropped" not "clipped") ---
$ws.Range("C12").Value = "sagittal"
$ws.Range("F12").Value = "_neck_cropped.obj"

# --- Row 13: a brand-new "patient" execution-parameter row, mostly a
#     copy of row 12 but pointing at the sample_data folder and with
#     export_landmarks turned on ---
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "patient"
$ws.Range("C13").Value = "sagittal"
$ws.Range("D13").Value = $true
$ws.Range("E13").Value = "C:\Users\franz\Documents\work\projects\arp\data\patient_data\sample_data"
$ws.Range("F13").Value = "_neck_cropped.obj"
$ws.Range("G13").Value = "C:\Users\franz\Documents\work\projects\arp\data\patient_data\patient_information.xlsx"
$ws.Range("H13").Value = "C:\Users\franz\Documents\work\projects\arp\data\patient_data\patient_information.xlsx"
$ws.Range("I13").Value = $false
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = "automatic"
$ws.Range("L13").Value = "_"
$ws.Range("M13").Value = $true
$ws.Range("N13").Value = $true
$ws.Range("O13").Value = $true

# --- restore the cursor/selection to where the author left it ---
$ws.Range("L17").Select()
